$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.1044623333333333
$ws.Range("H2").Value = 0.313387
$ws.Range("I2").Value = 0.3334131260000681
$ws.Range("J2").Value = 0.3334131260000681
$ws.Range("M2").Value = 27.30988266666667
$ws.Range("N2").Value = 81.929648
$ws.Range("O2").Value = 0.2168690090390243
$ws.Range("P2").Value = 0.2168690090390242
$ws.Range("Q2").Value = 2.852854066419556
$ws.Range("R2").Value = 25.675686597776
$ws.Range("S2").Value = 0.0723069742362381
$ws.Range("T2").Value = 0.07230697423623808

# Row 3
$ws.Range("G3").Value = 0.1044623333333333
$ws.Range("H3").Value = 0.313387
$ws.Range("I3").Value = 0.3334131260000681
$ws.Range("J3").Value = 0.3334131260000681
$ws.Range("M3").Value = 20.371636
$ws.Range("N3").Value = 61.114908
$ws.Range("O3").Value = 0.1617720795708915
$ws.Range("P3").Value = 0.1617720795708915
$ws.Range("Q3").Value = 2.128068630377333
$ws.Range("R3").Value = 19.152617673396
$ws.Range("S3").Value = 0.0539369347492627
$ws.Range("T3").Value = 0.0539369347492627

# Row 4
$ws.Range("G4").Value = 0.1044623333333333
$ws.Range("H4").Value = 0.313387
$ws.Range("I4").Value = 0.3334131260000681
$ws.Range("J4").Value = 0.3334131260000681
$ws.Range("M4").Value = 42.52135866666666
$ws.Range("N4").Value = 127.564076
$ws.Range("O4").Value = 0.3376640254953711
$ws.Range("P4").Value = 0.3376640254953711
$ws.Range("Q4").Value = 4.441880342823556
$ws.Range("R4").Value = 39.976923085412
$ws.Range("S4").Value = 0.1125816182781783
$ws.Range("T4").Value = 0.1125816182781783

# Row 5
$ws.Range("G5").Value = 0.1044623333333333
$ws.Range("H5").Value = 0.313387
$ws.Range("I5").Value = 0.3334131260000681
$ws.Range("J5").Value = 0.3334131260000681
$ws.Range("M5").Value = 6.417394333333333
$ws.Range("N5").Value = 19.252183
$ws.Range("O5").Value = 0.05096081761571768
$ws.Range("P5").Value = 0.05096081761571768
$ws.Range("Q5").Value = 0.6703759859801112
$ws.Range("R5").Value = 6.033383873821
$ws.Range("S5").Value = 0.01699100550477577
$ws.Range("T5").Value = 0.01699100550477577

# Row 6
$ws.Range("G6").Value = 0.1044623333333333
$ws.Range("H6").Value = 0.313387
$ws.Range("I6").Value = 0.3334131260000681
$ws.Range("J6").Value = 0.3334131260000681
$ws.Range("M6").Value = 29.307738
$ws.Range("N6").Value = 87.923214
$ws.Range("O6").Value = 0.2327340682789955
$ws.Range("P6").Value = 0.2327340682789955
$ws.Range("Q6").Value = 3.061554696202
$ws.Range("R6").Value = 27.553992265818
$ws.Range("S6").Value = 0.07759659323161315
$ws.Range("T6").Value = 0.07759659323161315

# Row 7
$ws.Range("G7").Value = 0.2088496666666667
$ws.Range("H7").Value = 0.626549
$ws.Range("I7").Value = 0.6665868739999319
$ws.Range("J7").Value = 0.6665868739999319
$ws.Range("M7").Value = 27.30988266666667
$ws.Range("N7").Value = 81.929648
$ws.Range("O7").Value = 0.2168690090390243
$ws.Range("P7").Value = 0.2168690090390242
$ws.Range("Q7").Value = 5.703659891639112
$ws.Range("R7").Value = 51.332939024752
$ws.Range("S7").Value = 0.1445620348027861
$ws.Range("T7").Value = 0.1445620348027861

# Row 8
$ws.Range("G8").Value = 0.2088496666666667
$ws.Range("H8").Value = 0.626549
$ws.Range("I8").Value = 0.6665868739999319
$ws.Range("J8").Value = 0.6665868739999319
$ws.Range("M8").Value = 20.371636
$ws.Range("N8").Value = 61.114908
$ws.Range("O8").Value = 0.1617720795708915
$ws.Range("P8").Value = 0.1617720795708915
$ws.Range("Q8").Value = 4.254609388054667
$ws.Range("R8").Value = 38.291484492492
$ws.Range("S8").Value = 0.1078351448216288
$ws.Range("T8").Value = 0.1078351448216288

# Row 9
$ws.Range("G9").Value = 0.2088496666666667
$ws.Range("H9").Value = 0.626549
$ws.Range("I9").Value = 0.6665868739999319
$ws.Range("J9").Value = 0.6665868739999319
$ws.Range("M9").Value = 42.52135866666666
$ws.Range("N9").Value = 127.564076
$ws.Range("O9").Value = 0.3376640254953711
$ws.Range("P9").Value = 0.3376640254953711
$ws.Range("Q9").Value = 8.88057158374711
$ws.Range("R9").Value = 79.925144253724
$ws.Range("S9").Value = 0.2250824072171927
$ws.Range("T9").Value = 0.2250824072171927

# Row 10
$ws.Range("G10").Value = 0.2088496666666667
$ws.Range("H10").Value = 0.626549
$ws.Range("I10").Value = 0.6665868739999319
$ws.Range("J10").Value = 0.6665868739999319
$ws.Range("M10").Value = 6.417394333333333
$ws.Range("N10").Value = 19.252183
$ws.Range("O10").Value = 0.05096081761571768
$ws.Range("P10").Value = 0.05096081761571768
$ws.Range("Q10").Value = 1.340270667385222
$ws.Range("R10").Value = 12.062436006467
$ws.Range("S10").Value = 0.03396981211094191
$ws.Range("T10").Value = 0.03396981211094191

# Row 11
$ws.Range("G11").Value = 0.2088496666666667
$ws.Range("H11").Value = 0.626549
$ws.Range("I11").Value = 0.6665868739999319
$ws.Range("J11").Value = 0.6665868739999319
$ws.Range("M11").Value = 29.307738
$ws.Range("N11").Value = 87.923214
$ws.Range("O11").Value = 0.2327340682789955
$ws.Range("P11").Value = 0.2327340682789955
$ws.Range("Q11").Value = 6.120911312054001
$ws.Range("R11").Value = 55.08820180848601
$ws.Range("S11").Value = 0.1551374750473823
$ws.Range("T11").Value = 0.1551374750473823
